# Updates cryptocurrency price/volume figures in the active worksheet
# to reflect the latest scrape (rows 2-51, columns B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.448.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.719.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5309'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.73%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06712'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2671'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07697'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.515'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.953.24'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.716.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.54%  '

$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8229'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.387.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '223.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.677'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.053'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.699'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1211'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.265'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05385'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.83%  '

$ws.Range("E31").Value = '  -0.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.491'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.436'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.638'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.871'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9598'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.391'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5888'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.153.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01650'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.811'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.56%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8436'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.859.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4585'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.004'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.111'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05202'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.74%  '
